# "Update countries & provincias Spain" -- refresh the COVID country table
# on the "Pais" sheet: a handful of countries overtook their neighbour in
# the ranking (so the two rows swap which country they display) and the
# case-count columns (B..H) for the affected rows get the refreshed
# totals. The "last updated" banner in A1 also advances to the new time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 13:23"

# Row 16
$ws.Range("B16").Value = 508389
$ws.Range("C16").Value = 4108
$ws.Range("D16").Value = 411840
$ws.Range("E16").Value = 67479
$ws.Range("G16").Value = 254
$ws.Range("H16").Value = 29070

# Row 32
$ws.Range("B32").Value = 160461
$ws.Range("C32").Value = 3109
$ws.Range("D32").Value = 120515
$ws.Range("E32").Value = 34411
$ws.Range("G32").Value = 68
$ws.Range("H32").Value = 5535

# Row 37
$ws.Range("B37").Value = 128405
$ws.Range("C37").Value = 214
$ws.Range("D37").Value = 125373
$ws.Range("E37").Value = 2812

# Row 41 -> Nepal
$ws.Range("A41").Value = "Nepal"
$ws.Range("B41").Value = 115358
$ws.Range("C41").Value = 3556
$ws.Range("D41").Value = 78780
$ws.Range("E41").Value = 35915
$ws.Range("G41").Value = 18
$ws.Range("H41").Value = 663

# Row 42 -> Kuwait
$ws.Range("A42").Value = "Kuwait"
$ws.Range("B42").Value = 111893
$ws.Range("D42").Value = 103802
$ws.Range("E42").Value = 7427
$ws.Range("H42").Value = 664

# Row 44
$ws.Range("B44").Value = 108608
$ws.Range("C44").Value = 1315
$ws.Range("D44").Value = 100007
$ws.Range("E44").Value = 8153
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 448

# Row 45
$ws.Range("B45").Value = 107213
$ws.Range("C45").Value = 638
$ws.Range("D45").Value = 93557
$ws.Range("E45").Value = 12603
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 1053

# Row 58
$ws.Range("B58").Value = 65881
$ws.Range("C58").Value = 1445
$ws.Range("E58").Value = 14280
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 2101

# Row 70 -> Libia
$ws.Range("A70").Value = "Libia"
$ws.Range("B70").Value = 44985
$ws.Range("C70").Value = 1164
$ws.Range("D70").Value = 25007
$ws.Range("E70").Value = 19322
$ws.Range("G70").Value = 12
$ws.Range("H70").Value = 656

# Row 71 -> Estado de Palestina
$ws.Range("A71").Value = "Estado de Palestina"
$ws.Range("B71").Value = 44684
$ws.Range("D71").Value = 38228
$ws.Range("E71").Value = 6069
$ws.Range("H71").Value = 387

# Row 93 -> Malasia
$ws.Range("A93").Value = "Malasia"
$ws.Range("B93").Value = 16880
$ws.Range("C93").Value = 660
$ws.Range("D93").Value = 11372
$ws.Range("E93").Value = 5345
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 163

# Row 94 -> Madagascar
$ws.Range("A94").Value = "Madagascar"
$ws.Range("B94").Value = 16754
$ws.Range("C94").Value = 28
$ws.Range("D94").Value = 16124
$ws.Range("E94").Value = 393
$ws.Range("H94").Value = 237

# Row 98
$ws.Range("B98").Value = 15307
$ws.Range("C98").Value = 15
$ws.Range("D98").Value = 13508
$ws.Range("E98").Value = 1484

# Row 140 -> Malta
$ws.Range("A140").Value = "Malta"
$ws.Range("B140").Value = 3937
$ws.Range("C140").Value = 93
$ws.Range("D140").Value = 3012
$ws.Range("E140").Value = 881
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 44

# Row 141 -> Estonia
$ws.Range("A141").Value = "Estonia"
$ws.Range("B141").Value = 3908
$ws.Range("C141").Value = 25
$ws.Range("D141").Value = 3015
$ws.Range("E141").Value = 825
$ws.Range("H141").Value = 68

# Row 142 -> Somalia
$ws.Range("A142").Value = "Somalia"
$ws.Range("B142").Value = 3864
$ws.Range("D142").Value = 3089
$ws.Range("E142").Value = 676
$ws.Range("H142").Value = 99

# Row 179 -> Gibraltar
$ws.Range("A179").Value = "Gibraltar"
$ws.Range("B179").Value = 499
$ws.Range("C179").Value = 13
$ws.Range("D179").Value = 424
$ws.Range("E179").Value = 75
$ws.Range("H179").Value = 0

# Row 180 -> Comoras
$ws.Range("A180").Value = "Comoras"
$ws.Range("B180").Value = 495
$ws.Range("D180").Value = 475
$ws.Range("E180").Value = 13
$ws.Range("H180").Value = 7

# Row 216 -> Islas Malvinas
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

# Row 217 -> Montserrat
$ws.Range("A217").Value = "Montserrat"
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
